$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper re-ran and re-ordered a few same-day fixtures, plus a new
# match (row 56) was appended at the end of the sheet.

# --- Swap the two matches played on 2023-09-01/02 (rows 12 <-> 13) ---
$tmp12 = $ws.Range("F12:V12").Value()
$tmp13 = $ws.Range("F13:V13").Value()
$ws.Range("F12:V12").Value = $tmp13
$ws.Range("F13:V13").Value = $tmp12

# --- Swap the two matches played on 2023-09-03 (rows 15 <-> 16) ---
$tmp15 = $ws.Range("F15:V15").Value()
$tmp16 = $ws.Range("F16:V16").Value()
$ws.Range("F15:V15").Value = $tmp16
$ws.Range("F16:V16").Value = $tmp15

# --- Swap the two matches played on 2023-10-08 (rows 47 <-> 48) ---
$tmp47 = $ws.Range("F47:V47").Value()
$tmp48 = $ws.Range("F48:V48").Value()
$ws.Range("F47:V47").Value = $tmp48
$ws.Range("F48:V48").Value = $tmp47

# --- Append the new match as row 56 ---
$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "morocco"
$ws.Range("C56").Value = "botola-pro"
$ws.Range("D56").Value = "2023-2024"
$ws.Range("E56").Value = 45233.85416666666
$ws.Range("F56").Value = "Jeunesse Sportive Soualem"
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = "Chabab Mohammedia"
$ws.Range("I56").Value = 2
$ws.Range("J56").Value = 2.75
$ws.Range("K56").Value = "02/11/2023 08:42"
$ws.Range("L56").Value = 2.68
$ws.Range("M56").Value = "03/11/2023 20:22"
$ws.Range("N56").Value = 2.76
$ws.Range("O56").Value = "02/11/2023 08:42"
$ws.Range("P56").Value = 2.61
$ws.Range("Q56").Value = "03/11/2023 20:20"
$ws.Range("R56").Value = 2.65
$ws.Range("S56").Value = "02/11/2023 08:42"
$ws.Range("T56").Value = 3.17
$ws.Range("U56").Value = "03/11/2023 20:22"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/morocco/botola-pro/jeunesse-sportive-soualem-chabab-mohammedia/ETVXim75/"

# Match the formatting used by the rest of the table: column A (bordered,
# bold, centered "Indice" style) and column E (datetime number format).
$ws.Range("A55").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("E55").Copy()
$ws.Range("E56").PasteSpecial(-4122)
$excel.CutCopyMode = 0
